$wb = $excel.ActiveWorkbook

# --- 1. Add the new row to the existing "Python" sheet (Sheet3) ---
$pythonSheet = $wb.Worksheets.Item("Python")
$pythonSheet.Range("C20").Value = "https://www.udemy.com/course/100-days-of-code/"
$pythonSheet.Range("B20").Value = "100 Days of Code - The Complete Python Pro Bootcamp for 2021"
$pythonSheet.Range("B20").Select()

# --- 2. Add a new "Projects" worksheet after the last sheet ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$projects = $wb.Worksheets.Add($null, $lastSheet)
$projects.Name = "Projects"

$projects.Range("C2").Value = "https://www.udemy.com/course/50-projects-50-days/"
$projects.Range("B2").Value = "50 Projects In 50 Days - HTML, CSS & JavaScript"

$projects.Range("C4").Value = "https://www.udemy.com/course/web-projects-with-vanilla-javascript/"
$projects.Range("B4").Value = "20 Web Projects With Vanilla JavaScript"

$projects.Range("C6").Value = "https://www.udemy.com/course/electron-from-scratch/"
$projects.Range("B6").Value = "Electron From Scratch: Build Desktop Apps With JavaScript"

$projects.Range("B:B").ColumnWidth = 90.65
$projects.Range("C:C").ColumnWidth = 81.5

$projects.Range("B6").Select()
$projects.Activate()
